$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scale column D (value) values from rows 2-33 by 10000 (万元 -> 元 style rescale)
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 * 10000
    }
}
